$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format temporarily so numeric-looking strings
# (e.g. "27.031.29", "20.64") are stored as literal text, matching the
# original inlineStr cells instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.031.29"
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").Value = "1.863.59"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "306.32"
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.5096"
$ws.Range("E7").Value = "  -2.30%  "
$ws.Range("D8").Value = "0.3747"
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("D9").Value = "0.07150"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").Value = "0.8867"
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("D11").Value = "20.64"
$ws.Range("D12").Value = "0.07577"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").Value = "1.863.97"
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("D14").Value = "5.307"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "89.40"
$ws.Range("E15").Value = "  -2.93%  "
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "0.000008430"
$ws.Range("E17").Value = "  -3.19%  "
$ws.Range("D18").Value = "14.09"
$ws.Range("E18").Value = "  -2.67%  "
$ws.Range("D19").Value = "0.9990"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "27.085.15"
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").Value = "5.032"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").Value = "2.104.51"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").Value = "10.53"
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("D24").Value = "6.461"
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("D25").Value = "1.834"
$ws.Range("E25").Value = "  -2.27%  "
$ws.Range("D26").Value = "147.77"
$ws.Range("E26").Value = "  -3.77%  "
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").Value = "2.105"
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("D29").Value = "112.73"
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").Value = "4.669"
$ws.Range("E30").Value = "  -3.78%  "
$ws.Range("D31").Value = "4.699"
$ws.Range("E31").Value = "  -3.68%  "
$ws.Range("D32").Value = "0.09089"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").Value = "0.05140"
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("D34").Value = "3.059"
$ws.Range("E34").Value = "  -3.68%  "
$ws.Range("D35").Value = "1.154"
$ws.Range("E35").Value = "  -6.22%  "
$ws.Range("D36").Value = "0.7281"
$ws.Range("E36").Value = "  -5.29%  "
$ws.Range("D37").Value = "0.02039"
$ws.Range("E37").Value = "  -2.11%  "
$ws.Range("D38").Value = "2.493"
$ws.Range("E38").Value = "  -5.30%  "
$ws.Range("D39").Value = "3.032"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("D41").Value = "0.5308"
$ws.Range("E41").Value = "  -3.55%  "
$ws.Range("D42").Value = "6.577"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").Value = "116.39"
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("D44").Value = "8.272"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").Value = "0.1470"
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("D46").Value = "0.9988"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "0.4622"
$ws.Range("E47").Value = "  -3.50%  "
$ws.Range("D48").Value = "10.01"
$ws.Range("E48").Value = "  -3.47%  "
$ws.Range("D49").Value = "1.567"
$ws.Range("E49").Value = "  -2.73%  "
$ws.Range("D50").Value = "36.54"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "63.83"
$ws.Range("E51").Value = "  -4.11%  "

# Reset column D back to the default (Normal) style so no stray
# number-format style lingers on the cells themselves.
$ws.Range("D2:D51").Style = "Normal"

Write-Output "Applied cryptos update"
